# Scheduled market-data refresh: updates computed price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# sheets with freshly scraped values. Pure data update, no formula/layout
# changes.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H80").Value = 2602.0527
$ws.Range("I80").Value = 2089
$ws.Range("J80").Value = 2975.182
$ws.Range("K80").Value = 6267
$ws.Range("L80").Value = 8925.545999999998
$ws.Range("M80").Value = -5269
$ws.Range("N80").Value = -10921.546
$ws.Range("H83").Value = 2602.0527
$ws.Range("I83").Value = 2089
$ws.Range("J83").Value = 2975.182
$ws.Range("K83").Value = 18801
$ws.Range("L83").Value = 26776.638
$ws.Range("M83").Value = -13809
$ws.Range("N83").Value = -36760.638

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 62980.766
$ws.Range("I74").Value = 202424.8
$ws.Range("J74").Value = 4879.0835
$ws.Range("K74").Value = 202424.8
$ws.Range("L74").Value = 4879.0835
$ws.Range("M74").Value = -201550.8
$ws.Range("N74").Value = -6627.0835
$ws.Range("H77").Value = 62980.766
$ws.Range("I77").Value = 202424.8
$ws.Range("J77").Value = 4879.0835
$ws.Range("K77").Value = 1012124
$ws.Range("L77").Value = 24395.4175
$ws.Range("M77").Value = -1007756
$ws.Range("N77").Value = -33131.4175
$ws.Range("H88").Value = 2050
$ws.Range("I88").Value = 2066.6667
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 2066.6667
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1660.6667
$ws.Range("N88").Value = -2812
$ws.Range("H91").Value = 2050
$ws.Range("I91").Value = 2066.6667
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 2066.6667
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -662.6667000000002
$ws.Range("N91").Value = -4808
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1833.625
$ws.Range("I86").Value = 1614.2106
$ws.Range("K86").Value = 1614.2106
$ws.Range("M86").Value = -491.2106000000001
$ws.Range("H88").Value = 62418.375
$ws.Range("J88").Value = 62418.375
$ws.Range("L88").Value = 62418.375
$ws.Range("N88").Value = -63230.375
$ws.Range("H89").Value = 1833.625
$ws.Range("I89").Value = 1614.2106
$ws.Range("K89").Value = 8071.053000000001
$ws.Range("M89").Value = -2455.053000000001
$ws.Range("H91").Value = 62418.375
$ws.Range("J91").Value = 62418.375
$ws.Range("L91").Value = 62418.375
$ws.Range("N91").Value = -65226.375
$ws.Range("H94").Value = 2008.2858
$ws.Range("I94").Value = 1843.4
$ws.Range("J94").Value = 2099.889
$ws.Range("K94").Value = 1843.4
$ws.Range("L94").Value = 2099.889
$ws.Range("M94").Value = -1392.4
$ws.Range("N94").Value = -3001.889

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 11195.667
$ws.Range("I32").Value = 1005
$ws.Range("J32").Value = 31577
$ws.Range("K32").Value = 1005
$ws.Range("L32").Value = 31577
$ws.Range("M32").Value = -689
$ws.Range("N32").Value = -32209
$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 23000
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 23000
$ws.Range("M82").Value = -4639
$ws.Range("N82").Value = -23722
$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 23000
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 23000
$ws.Range("M85").Value = -3752
$ws.Range("N85").Value = -25496

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3626.2593
$ws.Range("I51").Value = 2074.75
$ws.Range("J51").Value = 3896.087
$ws.Range("K51").Value = 6224.25
$ws.Range("L51").Value = 11688.261
$ws.Range("M51").Value = -5764.25
$ws.Range("N51").Value = -12608.261

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2475.8
$ws.Range("I31").Value = 2475.8
$ws.Range("K31").Value = 2475.8
$ws.Range("M31").Value = -2183.8
$ws.Range("H37").Value = 2475.8
$ws.Range("I37").Value = 2475.8
$ws.Range("K37").Value = 2475.8
$ws.Range("M37").Value = -2198.8
$ws.Range("H102").Value = 1283.8667
$ws.Range("I102").Value = 1235
$ws.Range("J102").Value = 1381.6
$ws.Range("K102").Value = 1235
$ws.Range("L102").Value = 1381.6
$ws.Range("M102").Value = 387
$ws.Range("N102").Value = -4625.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1324.5834
$ws.Range("I7").Value = 1334.35
$ws.Range("K7").Value = 1334.35
$ws.Range("M7").Value = -1222.35
$ws.Range("H16").Value = 1060
$ws.Range("J16").Value = 1425
$ws.Range("L16").Value = 1425
$ws.Range("N16").Value = -1765
$ws.Range("H82").Value = 1071.4482
$ws.Range("I82").Value = 856
$ws.Range("J82").Value = 1153.5238
$ws.Range("K82").Value = 856
$ws.Range("L82").Value = 1153.5238
$ws.Range("M82").Value = -495
$ws.Range("N82").Value = -1875.5238
$ws.Range("H85").Value = 1071.4482
$ws.Range("I85").Value = 856
$ws.Range("J85").Value = 1153.5238
$ws.Range("K85").Value = 856
$ws.Range("L85").Value = 1153.5238
$ws.Range("M85").Value = 392
$ws.Range("N85").Value = -3649.5238
$ws.Range("H122").Value = 4505
$ws.Range("I122").Value = 5830.8335
$ws.Range("J122").Value = 3781.818
$ws.Range("K122").Value = 17492.5005
$ws.Range("L122").Value = 11345.454
$ws.Range("M122").Value = -15042.5005
$ws.Range("N122").Value = -16245.454
$ws.Range("H126").Value = 1324.5834
$ws.Range("I126").Value = 1334.35
$ws.Range("K126").Value = 4003.05
$ws.Range("M126").Value = -1533.05

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 19250
$ws.Range("J63").Value = 19250
$ws.Range("L63").Value = 19250
$ws.Range("N63").Value = -20498
$ws.Range("H66").Value = 19250
$ws.Range("J66").Value = 19250
$ws.Range("L66").Value = 57750
$ws.Range("N66").Value = -63990
$ws.Range("H132").Value = 1309628.5
$ws.Range("I132").Value = 3099139
$ws.Range("J132").Value = 1909.3846
$ws.Range("K132").Value = 9297417
$ws.Range("L132").Value = 5728.1538
$ws.Range("M132").Value = -9294887
$ws.Range("N132").Value = -10788.1538
